# Netlist.xlsx edit: replace the R1/R2/C2 (lumped RC) network with the new
# distributed-parameter SOC1/L1 network (Sparameters para distribuidos),
# and drop the now-unused component rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-style cells that need to pick up an existing centered style ---
# C2 used to have its own "vertical-center only" style; it now shares the
# common centered style already used by B2 (and most of the table).
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# E3 picks up the scientific-notation + centered style that E4 used to have.
$ws.Range("E4").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# F3/G3 pick up the plain centered style already used by A3/B3/etc.
$ws.Range("A3").Copy()
$ws.Range("F3:G3").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row 2: SOC1 between N1 and N2 ---
$ws.Range("A2").Value = "SOC1"
$ws.Range("B2").Value = "N1"
$ws.Range("C2").Value = "N2"
$ws.Range("D2").Value = "SOC"
$ws.Range("E2").Value = 13.2629
$ws.Range("F2").Value = 45
$ws.Range("G2").Value = 1000000000

# --- Row 3: L1 between N2 and node "0" ---
$ws.Range("A3").Value = "L1"
$ws.Range("B3").Value = "N2"
$ws.Range("C3").Value = "0"
$ws.Range("D3").Value = "L"
$ws.Range("E3").Value = 0.00000001
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

# --- Rows 4 and 5: clear the now-unused components (old C2/R2 rows) ---
# A:E keep their formatting but lose their values; F:G are fully cleared
# (they carried no special formatting to begin with).
$ws.Range("A4:E5").ClearContents()
$ws.Range("F4:G5").Clear()

# Selection now covers the header + data block instead of the old stray H17.
$ws.Range("A1:G3").Select()

$wb.Save()
